$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B width change (bestFit, wider) ---
$ws.Range("B1").EntireColumn.ColumnWidth = 59.83203125

# --- Row 26: D:K values 1 -> 0 (style unchanged) ---
$ws.Range("D26:K26").Value = 0

# --- Rows with both value + style change ---
$rows = @(112, 113, 114, 129, 130, 131, 141, 142)
foreach ($r in $rows) {
    $ws.Range("D$r:K$r").Value = 0

    $rng = $ws.Range("C$r:K$r")
    $rng.Borders.LineStyle = 1
    $rng.Borders.Color = 11184810
    $rng.Borders(7).LineStyle = -4142
    $rng.Interior.Color = 16777215
    $rng.Interior.PatternColor = 0
    $rng.Font.Color = 0
}
